$d = $word.ActiveDocument

# Edit 1: "more abstract" -> "with more difficulty"
$d.Content.Find.Execute("more abstract", $true, $false, $false, $false, $false,
                         $true, 1, $false, "with more difficulty", 2)

# Edit 2: insert "core " before "concepts. These abstractions are invoked"
$d.Content.Find.Execute("re-usability of parallel concepts. These", $true, $false, $false, $false, $false,
                         $true, 1, $false, "re-usability of parallel core concepts. These", 2)

# Edit 3: "claims are unfounded" -> "claims were unfounded"
$d.Content.Find.Execute("claims are unfounded", $true, $false, $false, $false, $false,
                         $true, 1, $false, "claims were unfounded", 2)
